# Model version-compatibility bugfix
# - iterative calculation enabled (calcPr iterateDelta)
# - comment author/source correction on Europe!A15
# - refreshed historic CO2 emissions series (Europe!C15:AB15) + trimmed
#   trailing columns that no longer belong to the named range
# - view/selection bookkeeping (active tab, per-sheet selections)

$wb = $excel.ActiveWorkbook

$wsGlobal    = $wb.Worksheets.Item("Global")
$wsEurope    = $wb.Worksheets.Item("Europe")
$wsCatalonia = $wb.Worksheets.Item("Catalonia")

# --- Workbook calculation options -----------------------------------------
$excel.Iteration    = $true
$excel.MaxIterations = 100
$excel.MaxChange     = 0.0001

# --- Comment correction on Europe!A15 --------------------------------------
$cmt = $wsEurope.Range("A15").Comment
$newComment = "Iñigo:" + [char]10 + "CAIT" + [char]10 + "http://www.wri.org/resources/data-sets/cait-historical-emissions-data-countries-us-states-unfccc" + [char]10
$cmt.Text($newComment) | Out-Null

# --- Historic CO2 emissions from land use change and forestry (Europe) -----
$wsEurope.Range("C15").Value = -0.29501957469948298
$wsEurope.Range("D15").Value = -0.29501419551768499
$wsEurope.Range("E15").Value = -0.285226655822994
$wsEurope.Range("F15").Value = -0.29510965510000098
$wsEurope.Range("G15").Value = -0.29510446800000001
$wsEurope.Range("H15").Value = -0.29509928569999999
$wsEurope.Range("I15").Value = -0.29509653949999998
$wsEurope.Range("J15").Value = -0.29520885149999798
$wsEurope.Range("K15").Value = -0.29503016589999997
$wsEurope.Range("L15").Value = -0.29519845569999997
$wsEurope.Range("M15").Value = -0.29508998359999999
$wsEurope.Range("N15").Value = -0.32775546849999998
$wsEurope.Range("O15").Value = -0.32771818409999998
$wsEurope.Range("P15").Value = -0.32725808360000003
$wsEurope.Range("Q15").Value = -0.32795720229999997
$wsEurope.Range("R15").Value = -0.32785511909999898
$wsEurope.Range("S15").Value = -0.342626737499999
$wsEurope.Range("T15").Value = -0.34198419730000101
$wsEurope.Range("U15").Value = -0.342547157800001
$wsEurope.Range("V15").Value = -0.34320179969999998
$wsEurope.Range("W15").Value = -0.34296992100000101
$wsEurope.Range("X15").Value = -0.42800970090000101
$wsEurope.Range("Y15").Value = -0.42809877999999901
$wsEurope.Range("Z15").Value = -0.42900656040000001
$wsEurope.Range("AA15").Value = -0.42908028699999901

# AB15 keeps its number formatting but no longer holds a value, and the
# remaining AC15:AJ15 tail (which used to carry extra, unreferenced years)
# is dropped entirely.
$wsEurope.Range("AB15").ClearContents()
$wsEurope.Range("AC15:AJ15").ClearContents()

# --- Selections / active views ---------------------------------------------
$wsGlobal.Activate() | Out-Null
$wsGlobal.Range("C16").Select() | Out-Null

$wsEurope.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$wsEurope.Range("AB5").Select() | Out-Null

# Catalonia is activated last so it ends up as the workbook's active tab.
$wsCatalonia.Activate() | Out-Null
$wsCatalonia.Range("D3").Select() | Out-Null
